$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix bug update rank KH
# The shared text "220000" (referenced by D4:D7 and F16) becomes "8536362".
# Route the new value through Formula -> Copy -> PasteSpecial(values) so the
# result lands back as a literal shared-string cell (not a live formula) and
# the cell keeps its original style (a direct numeric-looking .Value=
# assignment would otherwise silently coerce the text into a Number cell).
foreach ($addr in @("D4", "D5", "D6", "D7", "F16")) {
    $cell = $ws.Range($addr)
    $cell.Formula = '="8536362"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

# Update the ranking threshold columns for rows 4-7: E 1 -> 13, F 1 -> 2
foreach ($row in 4..7) {
    $ws.Cells.Item($row, 5).Value = 13   # column E
    $ws.Cells.Item($row, 6).Value = 2    # column F
}
